$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the mesh file pointer from the 0.02 resolution mesh to the 0.05 one
$ws.Range("B8").Value = "data/new_area/mesh_0.05.msh"

# Move the active selection to B9 (as last interacted cell)
$ws.Range("B9").Select()
